# Apply "added inconsistent issues identified on Jan 11, 2022" edit.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # termsWithMulitpleLabels
$ws2 = $wb.Worksheets.Item(2)   # LabelsUsedMultipleTerms
$ws3 = $wb.Worksheets.Item(3)   # termWithDifferentParent
$ws4 = $wb.Worksheets.Item(4)   # units_issues

# ---------------------------------------------------------------------------
# 1) termsWithMulitpleLabels: append the new "multiple labels" rows (144-154)
#    Fill order (left-to-right, top-to-bottom) mirrors the original authoring
#    order so new shared-string indices land where the target file expects.
# ---------------------------------------------------------------------------
$ws1.Cells.Item(144,1).Value = "EUPATH_0000387"
$ws1.Cells.Item(144,2).Value = 2
$ws1.Cells.Item(144,3).Value = "Traveled in last month | Trips in last month"
$ws1.Cells.Item(144,4).Value = "icemr_south_asia | icemr_amazonia_peru"

$ws1.Cells.Item(145,1).Value = "EUPATH_0000407"
$ws1.Cells.Item(145,2).Value = 2
$ws1.Cells.Item(145,3).Value = "District in India | District in Zimbabwe"
$ws1.Cells.Item(145,4).Value = "icemr_india_cx | icemr_india_cohort | gates_shine"

$ws1.Cells.Item(146,1).Value = "EUPATH_0031266"
$ws1.Cells.Item(146,2).Value = 2
$ws1.Cells.Item(146,3).Value = "Rotavirus vaccine | Rotavirus vaccine_undocumented"
$ws1.Cells.Item(146,4).Value = "gates_vida | gates_vida_hucs_kenya | gates_vida_hucs_gambia_mali | gates_avenir"

$ws1.Cells.Item(147,1).Value = "EUPATH_0033175"
$ws1.Cells.Item(147,2).Value = 2
$ws1.Cells.Item(147,3).Value = "Diphtheria, pertussis, and tetanus (DPT) vaccine status | Diphtheria, pertussis, and tetanus (DPT) vaccine status, undocumented"
$ws1.Cells.Item(147,4).Value = "gates_perch | gates_avenir"

$ws1.Cells.Item(148,1).Value = "EUPATH_0033225"
$ws1.Cells.Item(148,2).Value = 2
$ws1.Cells.Item(148,3).Value = "Pneumococcal conjugate vaccine (PCV) status | Pneumococcal conjugate vaccine (PCV) status, undocumented"
$ws1.Cells.Item(148,4).Value = "gates_perch | gates_avenir"

$ws1.Cells.Item(149,1).Value = "EUPATH_0036100"
$ws1.Cells.Item(149,2).Value = 3
$ws1.Cells.Item(149,3).Value = "Bacille Calmette-Guerin (BCG) vaccine | BCG vaccine administered | Bacille Calmette-Guerin (BCG) vaccine, undocumented"
$ws1.Cells.Item(149,4).Value = "gates_provide | gates_avenir"

$ws1.Cells.Item(150,1).Value = "EUPATH_0042044"
$ws1.Cells.Item(150,2).Value = 2
$ws1.Cells.Item(150,3).Value = "Syphilis test performed | Syphilis sample collected"
$ws1.Cells.Item(150,4).Value = "gates_pcs | general/general_promote"

$ws1.Cells.Item(151,1).Value = "EUPATH_0042153"
$ws1.Cells.Item(151,2).Value = 2
$ws1.Cells.Item(151,3).Value = "Treponema pallidum, by rapid test | Treponema p. pallidum, by RDT"
$ws1.Cells.Item(151,4).Value = "gates_pcs | general/general_promote"

$ws1.Cells.Item(152,1).Value = "EUPATH_0047251"
$ws1.Cells.Item(152,2).Value = 2
$ws1.Cells.Item(152,3).Value = "Baseline date | Baseline survey year"
$ws1.Cells.Item(152,4).Value = "gates_shine | gates_ganc"

$ws1.Cells.Item(153,1).Value = "EUPATH_0049211"
$ws1.Cells.Item(153,2).Value = 2
$ws1.Cells.Item(153,3).Value = "Verification of birth date | Verification of birth date,youth"
$ws1.Cells.Item(153,4).Value = "gates_gamin | gates_avenir"

$ws1.Cells.Item(154,1).Value = "EUPATH_0049212"
$ws1.Cells.Item(154,2).Value = 2
$ws1.Cells.Item(154,3).Value = "Person providing birth date information | Person providing birth date information, youth"
$ws1.Cells.Item(154,4).Value = "gates_gamin | gates_avenir"

# ---------------------------------------------------------------------------
# 2) termWithDifferentParent: append the new "different parent" rows (182-188)
# ---------------------------------------------------------------------------
$ws3.Cells.Item(182,1).Value = "OBI_0001169"
$ws3.Cells.Item(182,2).Value = 2
$ws3.Cells.Item(182,3).Value = "Age"
$ws3.Cells.Item(182,4).Value = "Observation details|Demographics"
$ws3.Cells.Item(182,5).Value = "gates_gamin | gates_gems_huas | gates_score_moz | icemr_prism2_border_cohort | gates_gems1a | general/general_umsp | gates_gems | gates_vida | gates_elicit | icemr_india_cx | icemr_amazonia_brazil | icemr_india_meghalaya | gates_washb_bangladesh | gates_score_nig | gates_gems1a_huas | gates_perch | general/general_promote | icemr_india_cohort | gates_score_burundi | gates_vida_hucs_kenya | gates_sip | gates_score_five_country | gates_score_seasonal | icemr_prism2 | gates_shine | icemr_india_behavior | icemr_prism | gates_provide | gates_score_sm_cohort | gates_llineup | gates_vida_hucs_gambia_mali | gates_score_rwanda | icemr_southern_africa | gates_avenir | gates_maled | gates_score_zanzibar | icemr_india_fever_surv"

$ws3.Cells.Item(183,1).Value = "EUPATH_0042399"
$ws3.Cells.Item(183,2).Value = 2
$ws3.Cells.Item(183,3).Value = "Child sex"
$ws3.Cells.Item(183,4).Value = "Obstetric history|Child demographics"
$ws3.Cells.Item(183,5).Value = "general/general_promote | gates_shine"

$ws3.Cells.Item(184,1).Value = "EUPATH_0010420"
$ws3.Cells.Item(184,2).Value = 2
$ws3.Cells.Item(184,3).Value = "Consent for enrollment given"
$ws3.Cells.Item(184,4).Value = "Eligibility criteria|Household administrative information"
$ws3.Cells.Item(184,5).Value = "gates_gamin | icemr_malawi | gates_avenir"

$ws3.Cells.Item(185,1).Value = "EUPATH_0042257"
$ws3.Cells.Item(185,2).Value = 2
$ws3.Cells.Item(185,3).Value = "Cough"
$ws3.Cells.Item(185,4).Value = "Child symptoms|Symptoms"
$ws3.Cells.Item(185,5).Value = "general/general_promote | gates_avenir"

$ws3.Cells.Item(186,1).Value = "EUPATH_0042260"
$ws3.Cells.Item(186,2).Value = 2
$ws3.Cells.Item(186,3).Value = "Diarrhea"
$ws3.Cells.Item(186,4).Value = "Child symptoms|Symptoms"
$ws3.Cells.Item(186,5).Value = "general/general_promote | gates_avenir"

$ws3.Cells.Item(187,1).Value = "EUPATH_0022088"
$ws3.Cells.Item(187,2).Value = 2
$ws3.Cells.Item(187,3).Value = "Hospitalization reason"
$ws3.Cells.Item(187,4).Value = "Clinical history|Hospitalization"
$ws3.Cells.Item(187,5).Value = "icemr_amazonia_brazil | gates_avenir"

$ws3.Cells.Item(188,1).Value = "EUPATH_0042044"
$ws3.Cells.Item(188,2).Value = 2
$ws3.Cells.Item(188,3).Value = "Syphilis test performed|Syphilis sample collected"
$ws3.Cells.Item(188,4).Value = "Blood bacteria testing|Blood sample collection process"
$ws3.Cells.Item(188,5).Value = "gates_pcs | general/general_promote"

# ---------------------------------------------------------------------------
# 3) Section header + re-stated column headers on termsWithMulitpleLabels
#    (row 132/135/140 style: bold divider line copied down, then header row)
# ---------------------------------------------------------------------------
$ws1.Range("A140").Copy()
$ws1.Range("A142").PasteSpecial(-4122)
$ws1.Cells.Item(142,1).Value = "Checked on Jan 11th, 2022, identified following issue:"

$ws1.Cells.Item(143,1).Value = "sid"
$ws1.Cells.Item(143,2).Value = "label_count"
$ws1.Cells.Item(143,3).Value = "labels"
$ws1.Cells.Item(143,4).Value = "datasets"

# ---------------------------------------------------------------------------
# 4) Section header + re-stated column headers on termWithDifferentParent
# ---------------------------------------------------------------------------
$ws3.Range("A177").Copy()
$ws3.Range("A180").PasteSpecial(-4122)
$ws3.Cells.Item(180,1).Value = "Checked on Jan 11th, 2022, identified following issue:"

$ws3.Cells.Item(181,1).Value = "sid"
$ws3.Cells.Item(181,2).Value = "pid_count"
$ws3.Cells.Item(181,3).Value = "labels"
$ws3.Cells.Item(181,4).Value = "pLabels"
$ws3.Cells.Item(181,5).Value = "datasets"

# ---------------------------------------------------------------------------
# 5) LabelsUsedMultipleTerms: "no issue was found" divider line
# ---------------------------------------------------------------------------
$ws2.Range("A69").Copy()
$ws2.Range("A71").PasteSpecial(-4122)
$ws2.Cells.Item(71,1).Value = "Checked on Jan 11th, 2022, no issue was found"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 6) Column C on termsWithMulitpleLabels needs to widen for the longer text
# ---------------------------------------------------------------------------
$ws1.Columns.Item(3).ColumnWidth = 74.65

# ---------------------------------------------------------------------------
# 7) Restore the active-sheet/selection state to match the saved workbook:
#    tab moves from units_issues to termWithDifferentParent.
# ---------------------------------------------------------------------------
$ws1.Range("A142").Select()
$ws2.Range("A75").Select()
$ws3.Activate()
$ws3.Range("C183").Select()
